# Add the "20-02-2023" consultant meeting entry at the end of the log.
#
# The document currently ends with a paragraph:
#   "Stored in: 13-02-2023-meeting.docx"
# which also carries the "_GoBack" bookmark at its very end.
#
# We need to turn that into three paragraphs:
#   1) "Stored in: 13-02-2023-meeting.docx"                              (unchanged text, bookmark removed)
#   2) "Consultant Meeting – 20th February 2023; 14:00 – 15:00"          (new Heading1, bold, "th" superscript)
#   3) "Stored in: 20-02-2023-meeting.docx"                              (new, carries the relocated _GoBack bookmark)

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Pull the "_GoBack" bookmark out of the last paragraph; it will be
#    re-created later, inside the new trailing "Stored in" paragraph.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2. Insert a new paragraph right after the existing last paragraph
#    (i.e. right after "...13-02-2023-meeting.docx").
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$insertAt = $lastPara.Range.End - 1   # just before the pilcrow
$d.Range($insertAt, $insertAt).InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newStart = $newRange.Start

# ------------------------------------------------------------------
# 3. Build the text for the new heading paragraph and the new
#    "Stored in" paragraph in one shot (joined with a carriage return
#    so Word splits them into two separate, plain paragraphs), then
#    apply formatting/style only where required.
# ------------------------------------------------------------------
$enDash = [char]0x2013

$headingPrefix = "Consultant Meeting " + $enDash + " 20"
$headingSup    = "th"
$headingSuffix = " February 2023; 14:00 " + $enDash + " 15:00"
$headingFull   = $headingPrefix + $headingSup + $headingSuffix

$storedPrefix = "Stored in"
$storedSuffix = ": 20-02-2023-meeting.docx"
$storedFull   = $storedPrefix + $storedSuffix

$newRange.Text = $headingFull + [char]13 + $storedFull

# Paragraph holding the heading text -> Heading1 style, bold run.
$headingPara = $newPara
$headingPara.Style = "Heading 1"
$headingPara.Range.Font.Bold = 1

# Make the "th" run superscript (splits the bold run into three runs).
$supStart = $newStart + $headingPrefix.Length
$supEnd = $supStart + $headingSup.Length
$d.Range($supStart, $supEnd).Font.Superscript = $true

# ------------------------------------------------------------------
# 4. Re-create the "_GoBack" bookmark in the new trailing paragraph,
#    right between "Stored in" and ": 20-02-2023-meeting.docx".
# ------------------------------------------------------------------
$storedPara = $d.Paragraphs.Last
$storedStart = $storedPara.Range.Start
$bmPos = $storedStart + $storedPrefix.Length
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))

Write-Output "Added 20-02-2023 consultant meeting entry; paragraph count = $($d.Paragraphs.Count)"
